$d = $word.ActiveDocument

# 1) "that is used" -> "that are used" (also removes the gramStart/gramEnd
#    proofErr markers that wrapped "is", since Find/Replace on plain text
#    re-creates the run without them)
$d.Content.Find.Execute("that is used", $true, $false, $false, $false, $false, $true, 1, $false, "that are used", 2)

# 2) Append a trailing space run after "Aim:"
$d.Content.Find.Execute("Aim:", $true, $false, $false, $false, $false, $true, 1, $false, "Aim: ", 2)

# 3) Add a new empty paragraph right after the "Assignments" paragraph
$r = $d.Content.Find.Execute("Assignments", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$p = $d.Paragraphs | Where-Object { $_.Range.Text -match "Assignments" } | Select-Object -First 1
$p.Range.InsertParagraphAfter()
